# Zero out the pandemic mobility/case columns (C:J) for every province
# EXCEPT the 5 provinces of interest (Central Java, East Java, Jakarta,
# South Sulawesi, West Java - rows 7, 10, 14, 27, 30), per the commit:
# "implemented changes to the map window and shapefile to focus on the
# 5 provinces of interest in Indonesia"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$keepRows = @(7, 10, 14, 27, 30)

for ($r = 2; $r -le 35; $r++) {
    if ($keepRows -contains $r) {
        continue
    }

    # Columns C and D carry a thousands-separator number style ("#,##0")
    # on some rows; once the values collapse to 0 the style reverts back
    # to the workbook default (no explicit style index).
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.NumberFormat() -eq "#,##0") {
        $cCell.ClearFormats() | Out-Null
    }
    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.NumberFormat() -eq "#,##0") {
        $dCell.ClearFormats() | Out-Null
    }

    $ws.Range($ws.Cells.Item($r, 3), $ws.Cells.Item($r, 10)).Value = 0
}

# Column A now needs room for the longest province name ("Jakarta Special
# Capital Region"); Excel auto-fit this to ~27.57 characters.
$ws.Columns.Item(1).ColumnWidth = 26.66666667

# The view has scrolled down so row 17 is at the top and the user
# highlighted the West Java block (C31:J35) as the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("C31:J35").Select() | Out-Null
